# "Generate Report for Archive"
#
# The handoff status text moved from "Ready for handoff" to "In Translation"
# on every sheet that surfaces it (the shared string is reused by the
# Overview sheet's per-locale status columns and by each locale sheet's
# own "Status" column). Updating the cell values here collapses back onto
# a single shared string, exactly like Excel would do when you retype a
# cell that happens to match another cell's text elsewhere in the book.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The status column is now narrower text ("In Translation" vs. "Ready for
# handoff"), so the report's columns were re-sized to fit the new content
# on every sheet where that column is shown.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
